$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name / title date
$ws.Name = "Through 2022-10-27"

# Update the "October (through 10-26)" label -> "October (through 10-27)"
$ws.Range("A11").Value = "October (through 10-27)"

# Update row 9 (September) I9 value
$ws.Range("I9").Value = 163

# Update row 11 (October) values
$ws.Range("B11").Value = 27
$ws.Range("D11").Value = 67
$ws.Range("E11").Value = 57
$ws.Range("F11").Value = 52
$ws.Range("G11").Value = 133
$ws.Range("H11").Value = 170
$ws.Range("I11").Value = 97

# Update row 12 (Total) values
$ws.Range("B12").Value = 253
$ws.Range("D12").Value = 694
$ws.Range("E12").Value = 605
$ws.Range("F12").Value = 474
$ws.Range("G12").Value = 1034
$ws.Range("H12").Value = 1417
$ws.Range("I12").Value = 1373
